# Update Week 17 target depth data for Rams OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet (row 3 = "R") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 461
$wsOff.Range("C3").Value = 324
$wsOff.Range("D3").Value = 120
$wsOff.Range("E3").Value = 61
$wsOff.Range("F3").Value = 11
$wsOff.Range("G3").Value = 7

# --- DEF sheet (row 3 = "R") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 502
$wsDef.Range("C3").Value = 373
$wsDef.Range("D3").Value = 112
$wsDef.Range("E3").Value = 55
$wsDef.Range("F3").Value = 9
$wsDef.Range("G3").Value = 6
